$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 363
$ws1.Range("F7").Value = 1180
$ws1.Range("F8").Value = 444
$ws1.Range("F9").Value = 7127
$ws1.Range("F10").Value = 83
$ws1.Range("F12").Value = 2048
$ws1.Range("F13").Value = 7994
$ws1.Range("F16").Value = 5513
$ws1.Range("F18").Value = 2425
$ws1.Range("F19").Value = 1028
$ws1.Range("F20").Value = 4563
$ws1.Range("F25").Value = 380
$ws1.Range("F26").Value = 259
$ws1.Range("F27").Value = 10
$ws1.Range("F28").Value = 2373
$ws1.Range("F31").Value = 82
$ws1.Range("F32").Value = 152
$ws1.Range("F33").Value = 586
$ws1.Range("F34").Value = 7
$ws1.Range("F36").Value = 1504
$ws1.Range("F39").Value = 2335
$ws1.Range("F40").Value = 2216

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 94
$ws2.Range("F3").Value = 79
$ws2.Range("F4").Value = 65
$ws2.Range("F5").Value = 7
$ws2.Range("F6").Value = 27
$ws2.Range("F7").Value = 26

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1280

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1280
$ws4.Range("F7").Value = 94
$ws4.Range("F8").Value = 363
$ws4.Range("F9").Value = 1180
$ws4.Range("F10").Value = 444
$ws4.Range("F11").Value = 7127
$ws4.Range("F12").Value = 83
$ws4.Range("F14").Value = 2048
$ws4.Range("F15").Value = 7994
$ws4.Range("F18").Value = 5513
$ws4.Range("F20").Value = 2425
$ws4.Range("F21").Value = 1028
$ws4.Range("F22").Value = 4563
$ws4.Range("F25").Value = 79
$ws4.Range("F27").Value = 65
$ws4.Range("F28").Value = 380
$ws4.Range("F29").Value = 10
$ws4.Range("F30").Value = 2373
$ws4.Range("F33").Value = 82
$ws4.Range("F34").Value = 152
$ws4.Range("F35").Value = 7
$ws4.Range("F36").Value = 586
$ws4.Range("F37").Value = 7
$ws4.Range("F39").Value = 27
$ws4.Range("F40").Value = 1504
$ws4.Range("F43").Value = 2335
$ws4.Range("F44").Value = 26
$ws4.Range("F45").Value = 2216
